$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.133320331573486
$ws.Range("B1").Value = 2.202712059020996
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.251590490341187
$ws.Range("E1").Value = 1.083869338035583
